# Update the ToDo list: the "Error handling" task moves up from row 6 to
# row 4 (ahead of "Replace XNAMATH..."), and its description gains an
# extra sentence. The two tasks it displaces ("Replace XNAMATH..." and
# "Refactoring...") each shift down one row, keeping their own text/estimate.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = "Error handling - go on a robustness run, set standards for future work.  Ensure anythign that can go wrong is handled."
$ws.Range("B4").Value = 21

$ws.Range("A5").Value = "Replace XNAMATH with a Rorn maths library"
$ws.Range("B5").Value = 14

$ws.Range("A6").Value = "Refactoring - we need consistency across the board"
$ws.Range("B6").Value = 21

# Column A needs to widen to fit the longer text (matches the saved
# best-fit width for the new, longer task description).
$ws.Columns.Item(1).ColumnWidth = 105.6666666666667

# Move the active selection to A5, matching the saved view state.
$ws.Range("A5").Select() | Out-Null
